$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired data (rows 2..19), matching the order/values described by the diff.
$data = @(
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Buddy Hield", "SG,SF", "Golden State Warriors"),
    @("Onyeka Okongwu", "PF,C", "Atlanta Hawks"),
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Miles Bridges", "SF,PF", "Charlotte Hornets"),
    @("Luka Doncic", "PG,SG", "Dallas Mavericks"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Donovan Clingan", "C", "Portland Trail Blazers"),
    @("Deandre Ayton", "C", "Portland Trail Blazers"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
